$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('E2').Value = '2026-02-14 18:48:23'
$ws.Range('I2').Value = '33.3 mm'
$ws.Range('O2').Value = '-1.0 °C'
$ws.Range('E3').Value = '2026-02-14 18:48:25'
$ws.Range('I3').Value = '13.6 mm'
$ws.Range('N3').Value = '-6.5 °C 18:20 TU'
$ws.Range('O3').Value = '-5.0 °C'
$ws.Range('E4').Value = '2026-02-14 18:48:27'
$ws.Range('J4').Value = '995.9 hPa'
$ws.Range('E5').Value = '2026-02-14 18:48:34'
$ws.Range('I5').Value = '20.6 mm'
$ws.Range('N5').Value = '-6.1 °C 18:15 TU'
$ws.Range('E6').Value = '2026-02-14 18:48:37'
$ws.Range('H6').Value = '79%'
$ws.Range('J6').Value = '996.0 hPa'
$ws.Range('E7').Value = '2026-02-14 18:48:39'
$ws.Range('J7').Value = '996.1 hPa'
$ws.Range('O7').Value = '13.1 °C'
$ws.Range('E8').Value = '2026-02-14 18:48:42'
$ws.Range('H8').Value = '65%'
$ws.Range('J8').Value = '995.9 hPa'
$ws.Range('O8').Value = '9.7 °C'
$ws.Range('E9').Value = '2026-02-14 18:48:44'
$ws.Range('H9').Value = '55%'
$ws.Range('E10').Value = '2026-02-14 18:48:47'
$ws.Range('H10').Value = '80%'
$ws.Range('E11').Value = '2026-02-14 18:48:49'
$ws.Range('H11').Value = '61%'
$ws.Range('E12').Value = '2026-02-14 18:48:52'
$ws.Range('E13').Value = '2026-02-14 18:48:54'
$ws.Range('J13').Value = '998.4 hPa'
$ws.Range('E14').Value = '2026-02-14 18:48:57'
$ws.Range('E15').Value = '2026-02-14 18:48:59'
$ws.Range('E16').Value = '2026-02-14 18:49:02'
$ws.Range('E17').Value = '2026-02-14 18:49:04'
$ws.Range('H17').Value = '69%'
$ws.Range('O17').Value = '1.9 °C'
$ws.Range('E18').Value = '2026-02-14 18:49:07'
$ws.Range('H18').Value = '77%'
$ws.Range('J18').Value = '996.2 hPa'
$ws.Range('E19').Value = '2026-02-14 18:49:09'
$ws.Range('E20').Value = '2026-02-14 18:49:12'
$ws.Range('I20').Value = '1.5 mm'
$ws.Range('L20').Value = '106.9 km/h - 329º 18:27 TU'
$ws.Range('N20').Value = '-6.9 °C 18:21 TU'
$ws.Range('E21').Value = '2026-02-14 18:49:14'
$ws.Range('J21').Value = '998.3 hPa'
$ws.Range('L21').Value = '46.1 km/h - 26º 18:03 TU'
$ws.Range('E22').Value = '2026-02-14 18:49:16'
$ws.Range('N22').Value = '-8.6 °C 18:05 TU'
$ws.Range('E23').Value = '2026-02-14 18:49:19'
$ws.Range('I23').Value = '36.2 mm'
$ws.Range('N23').Value = '-7.8 °C 18:27 TU'
$ws.Range('E24').Value = '2026-02-14 18:49:21'
$ws.Range('H24').Value = '67%'
$ws.Range('J24').Value = '1000.0 hPa'
$ws.Range('E25').Value = '2026-02-14 18:49:24'
$ws.Range('I25').Value = '11.5 mm'
$ws.Range('N25').Value = '-6.5 °C 18:14 TU'
$ws.Range('O25').Value = '-4.5 °C'
$ws.Range('E26').Value = '2026-02-14 18:49:26'
$ws.Range('E27').Value = '2026-02-14 18:49:29'
$ws.Range('N27').Value = '-4.2 °C 18:25 TU'
$ws.Range('E28').Value = '2026-02-14 18:49:32'
$ws.Range('H28').Value = '71%'
$ws.Range('J28').Value = '995.9 hPa'
$ws.Range('E29').Value = '2026-02-14 18:49:34'
$ws.Range('E30').Value = '2026-02-14 18:49:37'
$ws.Range('J30').Value = '995.8 hPa'
$ws.Range('E31').Value = '2026-02-14 18:49:39'
$ws.Range('J31').Value = '994.9 hPa'
$ws.Range('E32').Value = '2026-02-14 18:49:41'
$ws.Range('K32').Value = '11.0 MJ/m2'
$ws.Range('E33').Value = '2026-02-14 18:49:44'
$ws.Range('J33').Value = '997.7 hPa'
$ws.Range('N33').Value = '1.4 °C 18:29 TU'
$ws.Range('O33').Value = '4.2 °C'
$ws.Range('E34').Value = '2026-02-14 18:49:47'
$ws.Range('N34').Value = '-4.0 °C 18:28 TU'
$ws.Range('E35').Value = '2026-02-14 18:49:49'
$ws.Range('J35').Value = '1002.8 hPa'
$ws.Range('E36').Value = '2026-02-14 18:49:52'
$ws.Range('J36').Value = '996.6 hPa'
$ws.Range('E37').Value = '2026-02-14 18:49:54'
$ws.Range('H37').Value = '68%'
$ws.Range('J37').Value = '996.8 hPa'
$ws.Range('O37').Value = '6.8 °C'
$ws.Range('E38').Value = '2026-02-14 18:49:57'
$ws.Range('E39').Value = '2026-02-14 18:49:59'
$ws.Range('I39').Value = '11.8 mm'
$ws.Range('N39').Value = '-7.6 °C 18:29 TU'
$ws.Range('O39').Value = '-5.5 °C'
$ws.Range('E40').Value = '2026-02-14 18:50:02'
$ws.Range('I40').Value = '0.5 mm'
$ws.Range('J40').Value = '998.8 hPa'
$ws.Range('E41').Value = '2026-02-14 18:50:04'
$ws.Range('J41').Value = '997.9 hPa'
$ws.Range('O41').Value = '13.4 °C'
$ws.Range('E42').Value = '2026-02-14 18:50:07'
$ws.Range('E43').Value = '2026-02-14 18:50:09'
$ws.Range('E44').Value = '2026-02-14 18:50:12'
$ws.Range('I44').Value = '35.3 mm'
$ws.Range('N44').Value = '-6.5 °C 18:20 TU'
$ws.Range('E45').Value = '2026-02-14 18:50:14'
$ws.Range('J45').Value = '1004.9 hPa'
$ws.Range('E46').Value = '2026-02-14 18:50:17'
$ws.Range('J46').Value = '1001.3 hPa'
